$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the original inline-string cells).
$textForcedCells = @{
  "D5" = "229.25"
  "D7" = "60.95"
  "D9" = "0.386"
  "D10" = "0.0804"
  "D12" = "14.77"
  "D14" = "21.09"
  "D20" = "69.84"
  "D22" = "225.95"
  "D25" = "2.23"
  "D26" = "165.85"
  "D29" = "19.04"
  "D32" = "4.53"
  "D36" = "6.38"
  "D41" = "17.05"
  "D42" = "97.64"
  "D45" = "0.0922"
  "D47" = "4.05"
  "D50" = "7.01"
}
foreach ($addr in $textForcedCells.Keys) {
  $ws.Range($addr).NumberFormat = "@"
  $ws.Range($addr).Value = $textForcedCells[$addr]
  $ws.Range($addr).Style = "Normal"
}

# Remaining cell updates (safe as plain text, no numeric coercion risk).
$cellValues = @{
  "D2" = "38.134.77"
  "D3" = "2.054.97"
  "E3" = "  +1.49%  "
  "E4" = "  +0.08%  "
  "E5" = "  +0.29%  "
  "E6" = "  +0.72%  "
  "E7" = "  +8.68%  "
  "E8" = "  -0.06%  "
  "E10" = "  +2.98%  "
  "E11" = "  +1.67%  "
  "E12" = "  +3.39%  "
  "D13" = "2.357.90"
  "E13" = "  +1.53%  "
  "E14" = "  +5.13%  "
  "E15" = "  +2.73%  "
  "E16" = "  +2.49%  "
  "D17" = "2.048.38"
  "E17" = "  +1.21%  "
  "D18" = "38.060.67"
  "E18" = "  +2.38%  "
  "E19" = "  +2.48%  "
  "E20" = "  +1.19%  "
  "D21" = "0.0₃0832"
  "E21" = "  +1.85%  "
  "E22" = "  +1.28%  "
  "E23" = "  +0.00%  "
  "E24" = "  +0.00%  "
  "E25" = "  +0.61%  "
  "E26" = "  +1.44%  "
  "E27" = "  +2.04%  "
  "E28" = "  +4.02%  "
  "E29" = "  +1.75%  "
  "E30" = "  -1.24%  "
  "E31" = "  +2.11%  "
  "E32" = "  +1.55%  "
  "E33" = "  +2.84%  "
  "E34" = "  +8.35%  "
  "E35" = "  +0.71%  "
  "E36" = "  +16.20%  "
  "E37" = "  -2.94%  "
  "E38" = "  +3.02%  "
  "E39" = "  +0.17%  "
  "D40" = "1.518.17"
  "E40" = "  +3.03%  "
  "E41" = "  +4.95%  "
  "E42" = "  +3.15%  "
  "E43" = "  +1.16%  "
  "E44" = "  +2.41%  "
  "E45" = "  +0.77%  "
  "E46" = "  +1.52%  "
  "E47" = "  -3.59%  "
  "E48" = "  +0.80%  "
  "E49" = "  +1.39%  "
  "E50" = "  -1.63%  "
  "D51" = "2.246.88"
  "E51" = "  +1.74%  "
}
foreach ($addr in $cellValues.Keys) {
  $ws.Range($addr).Value = $cellValues[$addr]
}
